# Update cryptos list data (prices + 1h volume change) per latest scrape.
$wbx = $excel.ActiveWorkbook
$ws = $wbx.ActiveSheet

$ws.Range("D2").Value = "42.275.45"
$ws.Range("E2").Value = "  -3.44%  "
$ws.Range("D3").Value = "2.240.74"
$ws.Range("E3").Value = "  -4.63%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.13"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.08%  "
$ws.Range("E6").Value = "  -5.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.02"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.57%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.567"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -4.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0994"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.90"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "35.58"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +8.42%  "
$ws.Range("E13").Value = "  -2.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.81"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -7.00%  "
$ws.Range("D15").Value = "2.577.99"
$ws.Range("E15").Value = "  -4.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.95"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -8.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.869"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.60%  "
$ws.Range("D18").Value = "2.240.90"
$ws.Range("E18").Value = "  -4.55%  "
$ws.Range("D19").Value = "42.108.79"
$ws.Range("E19").Value = "  -3.67%  "
$ws.Range("D20").Value = "0.0₃0986"
$ws.Range("E20").Value = "  -2.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.93"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.22"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -7.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.31"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -7.19%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.90"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.67"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.05%  "
$ws.Range("E27").Value = "  -5.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.02"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -5.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.15"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.27"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.77"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -7.85%  "
$ws.Range("E32").Value = "  -6.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.127"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -6.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.39"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0718"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.81%  "
$ws.Range("E36").Value = "  -6.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.60"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.10"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +15.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.07"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.25"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -5.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0266"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "66.07"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.13"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.88"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.101"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -8.59%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.190"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.58%  "
$ws.Range("B47").Value = "BinanceUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("B48").Value = "SynthetixNetwork"
$ws.Range("C48").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.49"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +8.46%  "
$ws.Range("E49").Value = "  -4.67%  "
$ws.Range("B50").Value = "Celestia"
$ws.Range("C50").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.17"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +8.56%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.18"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.02%  "
